# Updated the Umple Game Project with initial observations
#
# - D14 header text corrected: "Technologies to used" -> "Technologies to be used"
# - Row 17 filled in with the third project's (Distributed Umple Game) initial
#   observations: requirements/feasibility/technologies/required-knowledge/SDLC type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Technologies to used" -> "Technologies to be used" header typo.
$ws.Range("D14").Value = "Technologies to be used"

# Row 17: Distributed Umple Game - initial observations (3rd project row).
$ws.Range("A17").Value = 3

$ws.Range("B17").Value = "Distributed Umple Game:-`nRequirements:`n- Demonstrate the features and the power of Model Driven Development through the use of Umple`n- Must be published on Github with an open source license`n- Must follow Agile methodology with User stories driving successive development sprints`n- Project Team members to maintain detailed logs in their respective project wikis of outcomes of meetings`n- Modify workbased on the boundaries of Umple use."
$ws.Range("B17").WrapText = $true
$ws.Range("B17").VerticalAlignment = -4160

$ws.Range("C17").Value = "Risks:`nDistributed Game to be decided upon. Additionally the architecture (preferably a Client Server implementation - in the interest of time) has to be decided upon`nFamiliarity with Umple.`nUmple to be used for code modification using native constructs. (This is not necessarily a risk with regards to the project  as it could lead to the potential discovery of bugs)`n"
$ws.Range("C17").WrapText = $true
$ws.Range("C17").VerticalAlignment = -4160

$ws.Range("D17").Value = "C++/ Java/ PHP or any of the Umple compliant languages - limited direct use.`nUmple should be used to the maximum possible extent for code development. Generated code can be in either or all of the above languages"
$ws.Range("D17").WrapText = $true
$ws.Range("D17").VerticalAlignment = -4160

$ws.Range("E17").Value = "Umple Engine"
$ws.Range("E17").VerticalAlignment = -4160

$ws.Range("F17").Value = "Agile"
$ws.Range("F17").VerticalAlignment = -4160

# Move the cursor / scroll the view to roughly where the edits were made.
$ws.Range("F17").Select()
